$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows (14-21): "Canje de Puntos" (Voucher + Credito) and "Resumen de Puntos" cases,
# following the same 3-value rotation (MIX/POS/PRE) used by the existing rows.
$data = @(
  @(14, "Canje_de_Puntos_Canje_de_Voucher_MIX", 1162816939),
  @(15, "Canje_de_Puntos_Canje_de_Voucher_POS", 1145642605),
  @(16, "Canje_de_Puntos_Canje_de_Voucher_PRE", 1162676705),
  @(17, "Canje_de_Puntos_Canje_de_Credito_MIX", 1162816939),
  @(18, "Canje_de_Puntos_Canje_de_Credito_PRE", 1162676705),
  @(19, "Resumen_de_Puntos_MIX", 1162816939),
  @(20, "Resumen_de_Puntos_POS", 1145642605),
  @(21, "Resumen_de_Puntos_MIX", 1162676705)
)

# Use the already-formatted B13 cell (label column header style) as the format
# template so the new B-column cells pick up the same numeric/alignment style
# as the rest of the B2:B13 data column.
$ws.Range("B13").Copy() | Out-Null

foreach ($row in $data) {
  $r = $row[0]
  $label = $row[1]
  $val = $row[2]
  $ws.Cells.Item($r, 1).Value = $label
  $ws.Cells.Item($r, 2).Value = $val
  $ws.Cells.Item($r, 2).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}

$excel.CutCopyMode = 0

# Move the active selection, matching the saved view state of the edited file.
$ws.Range("F15").Select() | Out-Null
